$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.852.88'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '1.544.04'
$ws.Range("E3").Value = '  -1.48%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").Value = '''206.09'
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("E6").Value = '  -0.96%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("E8").Value = '  -0.59%  '
$ws.Range("D9").Value = '''21.42'
$ws.Range("E9").Value = '  -2.67%  '
$ws.Range("E10").Value = '  -0.72%  '
$ws.Range("E11").Value = '  -1.14%  '
$ws.Range("D12").Value = '1.762.29'
$ws.Range("E12").Value = '  -1.58%  '
$ws.Range("D13").Value = '1.542.42'
$ws.Range("E13").Value = '  -1.16%  '
$ws.Range("D14").Value = '''3.69'
$ws.Range("E14").Value = '  -1.45%  '
$ws.Range("D15").Value = '''0.510'
$ws.Range("E15").Value = '  -1.03%  '
$ws.Range("D16").Value = '26.841.93'
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").Value = '''61.30'
$ws.Range("E17").Value = '  -0.34%  '
$ws.Range("D18").Value = '''215.31'
$ws.Range("E18").Value = '  +0.25%  '
$ws.Range("E19").Value = '  -2.45%  '
$ws.Range("E20").Value = '  +0.68%  '
$ws.Range("E21").Value = '  +0.25%  '
$ws.Range("D22").Value = '''4.02'
$ws.Range("E22").Value = '  -2.67%  '
$ws.Range("D23").Value = '''9.17'
$ws.Range("E23").Value = '  -1.56%  '
$ws.Range("E24").Value = '  -2.59%  '
$ws.Range("D25").Value = '''153.04'
$ws.Range("E25").Value = '  -0.53%  '
$ws.Range("E26").Value = '  -1.76%  '
$ws.Range("D27").Value = '''14.86'
$ws.Range("E27").Value = '  -0.84%  '
$ws.Range("E28").Value = '  +0.21%  '
$ws.Range("E29").Value = '  -0.67%  '
$ws.Range("E30").Value = '  -1.23%  '
$ws.Range("E31").Value = '  -1.91%  '
$ws.Range("E32").Value = '  +1.34%  '
$ws.Range("D33").Value = '1.369.39'
$ws.Range("E33").Value = '  -2.36%  '
$ws.Range("D34").Value = '''2.94'
$ws.Range("E34").Value = '  +0.53%  '
$ws.Range("E35").Value = '  -1.02%  '
$ws.Range("D36").Value = '''0.961'
$ws.Range("E36").Value = '  +2.44%  '
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("E38").Value = '  +1.32%  '
$ws.Range("D39").Value = '''0.521'
$ws.Range("E39").Value = '  -1.24%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '''0.807'
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '''5.75'
$ws.Range("E41").Value = '  +8.22%  '
$ws.Range("E42").Value = '  +0.24%  '
$ws.Range("D43").Value = '''0.990'
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("E44").Value = '  +1.62%  '
$ws.Range("D45").Value = '''63.25'
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("E46").Value = '  -3.79%  '
$ws.Range("D47").Value = '1.676.74'
$ws.Range("D48").Value = '''84.23'
$ws.Range("E48").Value = '  -2.27%  '
$ws.Range("E49").Value = '  +3.82%  '
$ws.Range("E50").Value = '  -1.43%  '
$ws.Range("E51").Value = '  +0.27%  '
